# Apply the target edit to the BBVA calculator values workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row 12: "Ticker/Index Name: 29/05/2025 14:15-14:30" block ---
$ws.Range("A12").Value = "Ticker/Index Name: 29/05/2025 14:15-14:30"
$ws.Range("B12").Value = 10.49249
$ws.Range("C12").Value = 16.088978999999998
$ws.Range("D12").Value = 12.423120000000001
$ws.Range("E12").Value = 6.8448929999999999
$ws.Range("F12").Value = 4.4332510000000003
$ws.Range("G12").Value = 4.535018
$ws.Range("H12").Value = 3.6621779999999999
$ws.Range("I12").Value = 4.4559819999999997
$ws.Range("J12").Value = 3.2815639999999999

# --- Sheet view: scroll position + active selection moved ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D8").Select()
